$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 17
$ws.Range("G5").Value = 6405.77
$ws.Range("B10").Value = 28375.41
$ws.Range("F56").Value = 6
$ws.Range("G56").Value = 1255.08
$ws.Range("F59").Value = 19
$ws.Range("G59").Value = 1560.66
$ws.Range("F64").Value = 117
$ws.Range("G64").Value = 9496.889999999999
$ws.Range("F68").Value = 49
$ws.Range("G68").Value = 5640.88
$ws.Range("F71").Value = 323
$ws.Range("G71").Value = 20575.1
$ws.Range("F75").Value = 3
$ws.Range("G75").Value = 1110.54
$ws.Range("F80").Value = 11
$ws.Range("G80").Value = 2706.77
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 0
$ws.Range("F84").Value = 34
$ws.Range("G84").Value = 3483.64
$ws.Range("B90").Value = 183083.37
$ws.Range("F115").Value = 209
$ws.Range("G115").Value = 20233.29
$ws.Range("B117").Value = 14091.17
$ws.Range("F145").Value = 494
$ws.Range("G145").Value = 3947.06
$ws.Range("B147").Value = 15177.5
$ws.Range("F149").Value = 234
$ws.Range("G149").Value = 15163.2
$ws.Range("F150").Value = 40
$ws.Range("G150").Value = 1859.6
$ws.Range("B156").Value = 32255.61
$ws.Range("F186").Value = 78
$ws.Range("G186").Value = 3351.66
$ws.Range("B192").Value = 48706
$ws.Range("E192").Value = 39.8
$ws.Range("F192").Value = -144
$ws.Range("G192").Value = -4795.2
$ws.Range("B193").Value = 64973
$ws.Range("E193").Value = 35.4
$ws.Range("F193").Value = 2
$ws.Range("G193").Value = 66.59999999999999
$ws.Range("F213").Value = 8
$ws.Range("G213").Value = 685.4400000000001
$ws.Range("B216").Value = 42183.59
$ws.Range("B219").Value = 63565
$ws.Range("E219").Value = 109.19
$ws.Range("F219").Value = 60
$ws.Range("G219").Value = 6162.6
$ws.Range("B220").Value = 61610
$ws.Range("E220").Value = 122.71
$ws.Range("F220").Value = -58
$ws.Range("G220").Value = -5957.18
$ws.Range("B232").Value = 55356
$ws.Range("E232").Value = 54.04
$ws.Range("F232").Value = -158
$ws.Range("G232").Value = -7527.12
$ws.Range("B233").Value = 63510
$ws.Range("E233").Value = 50.66
$ws.Range("F233").Value = 117
$ws.Range("G233").Value = 5573.88
$ws.Range("F251").Value = 1
$ws.Range("G251").Value = 244.75
$ws.Range("F255").Value = 569
$ws.Range("G255").Value = 97486.77
$ws.Range("B260").Value = 192507.83
$ws.Range("F282").Value = 2
$ws.Range("G282").Value = 107.4
$ws.Range("F303").Value = 34
$ws.Range("G303").Value = 7170.26
$ws.Range("B304").Value = 179932.48
$ws.Range("F320").Value = 55
$ws.Range("G320").Value = 3775.75
$ws.Range("B330").Value = 28625.47
$ws.Range("F345").Value = 57
$ws.Range("G345").Value = 3500.37
$ws.Range("B346").Value = 26370.01
$ws.Range("F360").Value = 0
$ws.Range("G360").Value = 0
$ws.Range("B361").Value = 0
$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9
$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29
$ws.Range("B375").Value = 45718
$ws.Range("E375").Value = 19.38
$ws.Range("F375").Value = -294
$ws.Range("G375").Value = -4768.68
$ws.Range("B376").Value = 64927
$ws.Range("E376").Value = 17.26
$ws.Range("F376").Value = 106
$ws.Range("G376").Value = 1719.32
$ws.Range("B385").Value = 65067
$ws.Range("E385").Value = 15.65
$ws.Range("F385").Value = 126
$ws.Range("G385").Value = 1855.98
$ws.Range("B386").Value = 53595
$ws.Range("E386").Value = 17.61
$ws.Range("F386").Value = -335
$ws.Range("G386").Value = -4934.55
$ws.Range("B463").Value = 64833
$ws.Range("E463").Value = 34.9
$ws.Range("F463").Value = 95
$ws.Range("G463").Value = 3118.85
$ws.Range("B464").Value = 60025
$ws.Range("E464").Value = 37.22
$ws.Range("F464").Value = -98
$ws.Range("G464").Value = -3217.34
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 108
$ws.Range("G473").Value = 3545.64
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("F485").Value = 16
$ws.Range("G485").Value = 2807.52
$ws.Range("B488").Value = 30745.52
$ws.Range("F509").Value = 222
$ws.Range("G509").Value = 17844.36
$ws.Range("B510").Value = 23976.82
$ws.Range("F555").Value = 21
$ws.Range("G555").Value = 1460.76
$ws.Range("B560").Value = 4962.92
$ws.Range("B572").Value = 65079
$ws.Range("F572").Value = 6
$ws.Range("G572").Value = 245.22
$ws.Range("B573").Value = 65362
$ws.Range("F573").Value = 20
$ws.Range("G573").Value = 817.4
$ws.Range("F577").Value = 63
$ws.Range("G577").Value = 2708.37
$ws.Range("F578").Value = 84
$ws.Range("G578").Value = 4190.76
$ws.Range("F582").Value = 37
$ws.Range("G582").Value = 2108.63
$ws.Range("B583").Value = 17599.08
$ws.Range("F599").Value = 1708
$ws.Range("G599").Value = 278591.88
$ws.Range("F601").Value = 414
$ws.Range("G601").Value = 117108.18
$ws.Range("B606").Value = 444861.21
$ws.Range("F613").Value = 140
$ws.Range("G613").Value = 22282.4
$ws.Range("B618").Value = 44379.11
$ws.Range("B619").Value = 1781869.81
$ws.Range("B620").Value = 1781869.81
